$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that currently sits at the end
#    of the "...debe marcar la opcion al iniciar." paragraph.
# ------------------------------------------------------------------
try {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
} catch {
    # no-op if it is not there
}

# ------------------------------------------------------------------
# 2) Locate the "Ud. acepta un empleo ..." paragraph (end of the
#    DISCUSION item) and append a brand-new paragraph right after it
#    holding the answer to the discussion question.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "Ud\. acepta un empleo") {
        $targetIndex = $i
    }
}

$target = $d.Paragraphs.Item($targetIndex)
$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newRange = $newPara.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="1416"/><w:jc w:val="both"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">Corregiría las ambigüedades en la interpretación de los requerimientos sin plantearlo como correcciones sino como mejoras, ya que existe una responsabilidad asumida con el actual empleador, y al solucionar dicho </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">inconveniente se lograría que la empresa o institución pueda alcanzar sus </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>objetivos</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>.</w:t></w:r></w:p>'

[void]$newRange.InsertXML($xml)

# ------------------------------------------------------------------
# 3) Footer page-count field: the extra paragraph pushes the doc from
#    7 to 8 pages, so refresh the cached PAGE field result shown in
#    the footer.
# ------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    $ftr = $sec.Footers.Item(1)
    $p = $ftr.Range.Paragraphs.Item(1)
    $chars = $p.Range.Characters
    for ($ci = 1; $ci -le $chars.Count; $ci++) {
        $ch = $chars.Item($ci)
        if ($ch.Text -eq "7") {
            $ch.Text = "8"
        }
    }
}

Write-Host "done"
